$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Reshape the "heats" worksheet from a wide layout (one column per
#    data point, one row per series: data/observation/dilution/deviation/
#    series) into a long layout (one row per data point, one column per
#    series: data/observation/dilution/deviation/series).
# ---------------------------------------------------------------------
$heats = $wb.Worksheets.Item("heats")

$usedRange = $heats.UsedRange
$lastCol = $usedRange.Columns.Count
$nPoints = $lastCol - 1

# Read the existing wide-format values (row 2..5 hold the series values,
# columns 2..lastCol hold one data point each).
$obs = @()
$dil = @()
$dev = @()
$series = @()
for ($c = 2; $c -le $lastCol; $c++) {
    $obs += , $heats.Cells.Item(2, $c).Value2
    $dil += , $heats.Cells.Item(3, $c).Value2
    $dev += , $heats.Cells.Item(4, $c).Value2
    $series += , $heats.Cells.Item(5, $c).Value2
}

# Wipe the old contents before writing the reshaped table.
$usedRange.ClearContents()

# New header row.
$heats.Cells.Item(1, 1).Value2 = "data"
$heats.Cells.Item(1, 2).Value2 = "observation"
$heats.Cells.Item(1, 3).Value2 = "dilution"
$heats.Cells.Item(1, 4).Value2 = "deviation"
$heats.Cells.Item(1, 5).Value2 = "series"

# New data rows: one row per original column/data point.
for ($i = 0; $i -lt $nPoints; $i++) {
    $r = $i + 2
    $heats.Cells.Item($r, 1).Value2 = $i + 1
    $heats.Cells.Item($r, 2).Value2 = $obs[$i]
    $heats.Cells.Item($r, 3).Value2 = $dil[$i]
    $heats.Cells.Item($r, 4).Value2 = $dev[$i]
    $heats.Cells.Item($r, 5).Value2 = $series[$i]
}

# ---------------------------------------------------------------------
# 2. Normalise each sheet's selection to a single cell (drops the stray
#    "4:4" row-range that had been tacked onto the sqref), and make the
#    "heats" sheet the active tab/selection.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("input_stoich_coefficients").Range("D4").Select() | Out-Null
$wb.Worksheets.Item("input_k_constants_log10").Range("A6").Select() | Out-Null
$wb.Worksheets.Item("input_concentrations").Range("E9").Select() | Out-Null
$wb.Worksheets.Item("setup").Range("A4").Select() | Out-Null
$wb.Worksheets.Item("enthalpies").Range("C8").Select() | Out-Null

$heats.Activate()
$heats.Range("C7").Select() | Out-Null
